$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.866.24'
$ws.Range("E2").Value = '  -1.09%  '

$ws.Range("D3").Value = '1.620.14'
$ws.Range("E3").Value = '  -1.77%  '

$ws.Range("E4").Value = '  -0.91%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.84'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.60%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.499'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.69%  '

$ws.Range("E7").Value = '  -0.63%  '

$ws.Range("E8").Value = '  -1.66%  '

$ws.Range("E9").Value = '  -3.37%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.45'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.29%  '

$ws.Range("E11").Value = '  -1.05%  '

$ws.Range("D12").Value = '1.845.32'
$ws.Range("E12").Value = '  -1.60%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.14'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.81%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.598.25'
$ws.Range("E14").Value = '  -3.99%  '

$ws.Range("E15").Value = '  -3.32%  '

$ws.Range("D16").Value = '25.879.55'
$ws.Range("E16").Value = '  -0.37%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.50'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.13%  '

$ws.Range("E18").Value = '  -3.28%  '

$ws.Range("E19").Value = '  -0.73%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '191.51'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.46%  '

$ws.Range("E21").Value = '  -2.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.49'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.88%  '

$ws.Range("E23").Value = '  -2.49%  '

$ws.Range("E24").Value = '  +2.68%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.32%  '

$ws.Range("E26").Value = '  -0.78%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.72'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.84%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.65'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.18%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.18'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.46%  '

$ws.Range("E30").Value = '  -1.51%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0476'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.62%  '

$ws.Range("E32").Value = '  -4.00%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.09'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.50%  '

$ws.Range("E34").Value = '  -2.71%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.48'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.08%  '

$ws.Range("D36").Value = '1.124.48'
$ws.Range("E36").Value = '  -0.77%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.839'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.93%  '

$ws.Range("E38").Value = '  -4.40%  '

$ws.Range("E39").Value = '  -2.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.510'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.30%  '

$ws.Range("E41").Value = '  -0.40%  '

$ws.Range("D42").Value = '1.755.75'
$ws.Range("E42").Value = '  -1.45%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.750'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.31%  '

$ws.Range("E44").Value = '  -5.56%  '

$ws.Range("E45").Value = '  -1.53%  '

$ws.Range("E46").Value = '  +1.39%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '53.96'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.29%  '

$ws.Range("E48").Value = '  -0.56%  '

$ws.Range("E49").Value = '  -1.75%  '

$ws.Range("E50").Value = '  -0.69%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.47'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.10%  '
